$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''66.276.29'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.24%  '
$ws.Range("D3").Value = '''3.030.59'
$ws.Range("D3").Style = "Normal"
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '''578.09'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.51%  '
$ws.Range("D6").Value = '''168.45'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.65%  '
$ws.Range("D7").Value = '''0.999'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("D8").Value = '''3.027.04'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.04%  '
$ws.Range("D9").Value = '''0.520'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.40%  '
$ws.Range("D10").Value = '''6.67'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.52%  '
$ws.Range("D11").Value = '''0.153'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.23%  '
$ws.Range("D12").Value = '''0.490'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +7.53%  '
$ws.Range("D13").Value = '''0.0000249'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.15%  '
$ws.Range("E14").Value = '  +6.70%  '
$ws.Range("D15").Value = '''0.125'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.05%  '
$ws.Range("D16").Value = '''66.271.56'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.26%  '
$ws.Range("D17").Value = '''3.532.57'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.15%  '
$ws.Range("D18").Value = '''7.23'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +4.68%  '
$ws.Range("D19").Value = '''16.48'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +19.07%  '
$ws.Range("D20").Value = '''3.027.47'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.07%  '
$ws.Range("D21").Value = '''468.51'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.39%  '
$ws.Range("D22").Value = '''0.710'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.54%  '
$ws.Range("D23").Value = '''7.39'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.61%  '
$ws.Range("D24").Value = '''83.06'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.89%  '
$ws.Range("D25").Value = '''12.76'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.19%  '
$ws.Range("E26").Value = '  -1.61%  '
$ws.Range("D27").Value = '''10.05'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.54%  '
$ws.Range("E28").Value = '  +0.04%  '
$ws.Range("D29").Value = '''8.21'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.35%  '
$ws.Range("D30").Value = '''2.42'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.72%  '
$ws.Range("D32").Value = '''0.118'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +7.20%  '
$ws.Range("D33").Value = '''0.0₃0998'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.17%  '
$ws.Range("E34").Value = '  +3.39%  '
$ws.Range("D35").Value = '''0.999'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.00%  '
$ws.Range("B36").Value = 'Filecoin'
$ws.Range("C36").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D36").Value = '''5.87'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.96%  '
$ws.Range("B37").Value = 'Mantle'
$ws.Range("C37").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D37").Value = '''0.992'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.13%  '
$ws.Range("D38").Value = '''48.20'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +9.12%  '
$ws.Range("D39").Value = '''2.07'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.62%  '
$ws.Range("D40").Value = '''49.49'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.13%  '
$ws.Range("D41").Value = '''0.312'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.67%  '
$ws.Range("E42").Value = '  -0.82%  '
$ws.Range("B43").Value = 'dogwifhat'
$ws.Range("C43").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D43").Value = '''2.85'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.71%  '
$ws.Range("B44").Value = 'Cosmos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D44").Value = '''8.65'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.32%  '
$ws.Range("E45").Value = '  +0.74%  '
$ws.Range("D46").Value = '''378.84'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.87%  '
$ws.Range("D47").Value = '''2.707.46'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.26%  '
$ws.Range("D48").Value = '''134.34'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.97%  '
$ws.Range("E49").Value = '  -0.01%  '
$ws.Range("D50").Value = '''24.51'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.95%  '
$ws.Range("E51").Value = '  +4.30%  '
